# Separate pre-mid and post-mid timetables:
# Regenerate the CSE semester-1 timetable workbook with the post-midsem
# room / section reassignment. Updates the Section A/B grids, the
# per-section course verification tables, the room-allocation summary
# (including the Auditorium split into C003 + room C404 being carved out
# of the old L408 lab block), the LTPSC compliance status glyphs, and the
# executive summary metrics.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Section_A timetable grid
# ---------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("B2").Value = "ELECTIVE_B1 [C104]"
$wsA.Range("C2").Value = "Free"
$wsA.Range("D2").Value = "ELECTIVE_B1 [C104]"
$wsA.Range("E2").Value = "MA162 [C304]"
$wsA.Range("F2").Value = "EC161 [C404]"

$wsA.Range("B3").Value = "MA162 [C304]"
$wsA.Range("C3").Value = "Free"
$wsA.Range("D3").Value = "HS161 [C101]"
$wsA.Range("E3").Value = "EC161 [C404]"
$wsA.Range("F3").Value = "MA161 [C002]"

$wsA.Range("B5").Value = "HS161 [C101]"
$wsA.Range("C5").Value = "Free"
$wsA.Range("D5").Value = "DS161 [C002]"
$wsA.Range("E5").Value = "MA161 [C002]"
$wsA.Range("F5").Value = "Free"

$wsA.Range("F6").Value = "Free"

$wsA.Range("C7").Value = "Free"
$wsA.Range("D7").Value = "EC161 (Lab) [L408]"
$wsA.Range("F7").Value = "DS161 [C002]"

$wsA.Range("D8").Value = "EC161 (Lab) [L408]"

# ---------------------------------------------------------------------
# Section_B timetable grid
# ---------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("B2").Value = "ELECTIVE_B1 [C002]"
$wsB.Range("C2").Value = "DS161 [C302]"
$wsB.Range("D2").Value = "ELECTIVE_B1 [C002]"
$wsB.Range("E2").Value = "Free"

$wsB.Range("B3").Value = "Free"
$wsB.Range("C3").Value = "EC161 [C003]"
$wsB.Range("D3").Value = "MA161 [C004]"
$wsB.Range("E3").Value = "MA161 [C004]"
$wsB.Range("F3").Value = "DS161 [C302]"

$wsB.Range("B5").Value = "MA162 [C401]"
$wsB.Range("C5").Value = "EC161 (Lab) [L408]"
$wsB.Range("D5").Value = "Free"
$wsB.Range("E5").Value = "Free"
$wsB.Range("F5").Value = "HS161 [C201]"

$wsB.Range("C6").Value = "EC161 (Lab) [L408]"

$wsB.Range("B7").Value = "Free"
$wsB.Range("C7").Value = "MA162 [C401]"
$wsB.Range("D7").Value = "Free"
$wsB.Range("E7").Value = "HS161 [C201]"
$wsB.Range("F7").Value = "EC161 [C003]"

$wsB.Range("E8").Value = "Free"

# ---------------------------------------------------------------------
# Verification_A course table
# ---------------------------------------------------------------------
$wsVA = $wb.Worksheets.Item("Verification_A")

$wsVA.Range("I2").Value = "C104"
$wsVA.Range("I3").Value = "C304"

$wsVA.Range("A4").Value = "**HS161**"
$wsVA.Range("B4").Value = "English Language"
$wsVA.Range("C4").Value = "Rajesh N S"
$wsVA.Range("D4").Value = "3-0-0-0-3"
$wsVA.Range("H4").Value = "Partial"
$wsVA.Range("I4").Value = "C101"

$wsVA.Range("A5").Value = "**DS161**"
$wsVA.Range("B5").Value = "Introduction to DATA science and artificial intelligence"
$wsVA.Range("C5").Value = "Girish Revadigar"
$wsVA.Range("D5").Value = "2-0-0-0-2"
$wsVA.Range("H5").Value = "Complete"
$wsVA.Range("I5").Value = "C002"

$wsVA.Range("A6").Value = "**EC161**"
$wsVA.Range("B6").Value = "Digital Design"
$wsVA.Range("C6").Value = "Prakash Pawar"
$wsVA.Range("F6").Value = "2/1"
$wsVA.Range("I6").Value = "C404, L408"

$wsVA.Range("A7").Value = "**MA161**"
$wsVA.Range("B7").Value = "Statistics"
$wsVA.Range("C7").Value = "Ramesh Athe"
$wsVA.Range("F7").Value = "0/0"
$wsVA.Range("I7").Value = "C002"

$wsVA.Range("H8").Value = "[WARN] 4 issues"

# ---------------------------------------------------------------------
# Verification_B course table
# ---------------------------------------------------------------------
$wsVB = $wb.Worksheets.Item("Verification_B")

$wsVB.Range("I2").Value = "C002"
$wsVB.Range("I3").Value = "C401"

$wsVB.Range("A4").Value = "**DS161**"
$wsVB.Range("B4").Value = "Introduction to DATA science and artificial intelligence"
$wsVB.Range("C4").Value = "Girish Revadigar"
$wsVB.Range("D4").Value = "2-0-0-0-2"
$wsVB.Range("H4").Value = "Complete"
$wsVB.Range("I4").Value = "C302"

$wsVB.Range("A5").Value = "**EC161**"
$wsVB.Range("B5").Value = "Digital Design"
$wsVB.Range("C5").Value = "Prakash Pawar"
$wsVB.Range("F5").Value = "2/1"
$wsVB.Range("I5").Value = "L408, C003"

$wsVB.Range("A6").Value = "**MA161**"
$wsVB.Range("B6").Value = "Statistics"
$wsVB.Range("C6").Value = "Ramesh Athe"
$wsVB.Range("D6").Value = "3-0-2-0-2"
$wsVB.Range("H6").Value = "Partial"
$wsVB.Range("I6").Value = "C004"

$wsVB.Range("A7").Value = "**HS161**"
$wsVB.Range("B7").Value = "English Language"
$wsVB.Range("C7").Value = "Rajesh N S"
$wsVB.Range("D7").Value = "3-0-0-0-3"
$wsVB.Range("F7").Value = "0/0"
$wsVB.Range("I7").Value = "C201"

$wsVB.Range("H8").Value = "[WARN] 4 issues"

# ---------------------------------------------------------------------
# Room_Allocation summary: insert a row so the Auditorium's large lecture
# block (C002) is split out, a new "large classroom" C003 appears, and the
# old L408 lab entry moves down to row 12 to make room for new row C404.
# ---------------------------------------------------------------------
$wsRA = $wb.Worksheets.Item("Room_Allocation")

$wsRA.Rows.Item(11).Insert()

$wsRA.Range("A2").Value = "C002"
$wsRA.Range("B2").Value = "large classroom"
$wsRA.Range("C2").Value = "116"
$wsRA.Range("D2").Value = "Projector"
$wsRA.Range("E2").Value = 6
$wsRA.Range("F2").Value = "A, B"
$wsRA.Range("G2").Value = 3
$wsRA.Range("H2").Value = "DS161, MA161, ELECTIVE_B1"
$wsRA.Range("I2").Value = "1.2"

$wsRA.Range("A3").Value = "C003"
$wsRA.Range("B3").Value = "large classroom"
$wsRA.Range("C3").Value = "135"
$wsRA.Range("D3").Value = "Projector"
$wsRA.Range("E3").Value = 2
$wsRA.Range("F3").Value = "B"
$wsRA.Range("G3").Value = 1
$wsRA.Range("H3").Value = "EC161"
$wsRA.Range("I3").Value = "0.4"

$wsRA.Range("A4").Value = "C004"
$wsRA.Range("B4").Value = "Auditorium"
$wsRA.Range("C4").Value = "240"
$wsRA.Range("D4").Value = "Projector"
$wsRA.Range("E4").Value = 2
$wsRA.Range("F4").Value = "B"
$wsRA.Range("G4").Value = 1
$wsRA.Range("H4").Value = "MA161"
$wsRA.Range("I4").Value = "0.4"

$wsRA.Range("A5").Value = "C101"
$wsRA.Range("B5").Value = "classroom"
$wsRA.Range("C5").Value = "96"
$wsRA.Range("D5").Value = "Projector"
$wsRA.Range("E5").Value = 2
$wsRA.Range("F5").Value = "A"
$wsRA.Range("G5").Value = 1
$wsRA.Range("H5").Value = "HS161"
$wsRA.Range("I5").Value = "0.4"

$wsRA.Range("A6").Value = "C104"
$wsRA.Range("B6").Value = "classroom"
$wsRA.Range("C6").Value = "96"
$wsRA.Range("D6").Value = "Projector"
$wsRA.Range("E6").Value = 2
$wsRA.Range("F6").Value = "A"
$wsRA.Range("G6").Value = 1
$wsRA.Range("H6").Value = "ELECTIVE_B1"
$wsRA.Range("I6").Value = "0.4"

$wsRA.Range("A7").Value = "C201"
$wsRA.Range("B7").Value = "classroom"
$wsRA.Range("C7").Value = "96"
$wsRA.Range("D7").Value = "Projector"
$wsRA.Range("E7").Value = 2
$wsRA.Range("F7").Value = "B"
$wsRA.Range("G7").Value = 1
$wsRA.Range("H7").Value = "HS161"
$wsRA.Range("I7").Value = "0.4"

$wsRA.Range("A8").Value = "C302"
$wsRA.Range("B8").Value = "classroom"
$wsRA.Range("C8").Value = "96"
$wsRA.Range("D8").Value = "Projector"
$wsRA.Range("E8").Value = 2
$wsRA.Range("F8").Value = "B"
$wsRA.Range("G8").Value = 1
$wsRA.Range("H8").Value = "DS161"
$wsRA.Range("I8").Value = "0.4"

$wsRA.Range("A9").Value = "C304"
$wsRA.Range("B9").Value = "classroom"
$wsRA.Range("C9").Value = "96"
$wsRA.Range("D9").Value = "Projector"
$wsRA.Range("E9").Value = 2
$wsRA.Range("F9").Value = "A"
$wsRA.Range("G9").Value = 1
$wsRA.Range("H9").Value = "MA162"
$wsRA.Range("I9").Value = "0.4"

$wsRA.Range("A10").Value = "C401"
$wsRA.Range("B10").Value = "classroom"
$wsRA.Range("C10").Value = "96"
$wsRA.Range("D10").Value = "Projector"
$wsRA.Range("E10").Value = 2
$wsRA.Range("F10").Value = "B"
$wsRA.Range("G10").Value = 1
$wsRA.Range("H10").Value = "MA162"
$wsRA.Range("I10").Value = "0.4"

$wsRA.Range("A11").Value = "C404"
$wsRA.Range("B11").Value = "classroom"
$wsRA.Range("C11").Value = "78"
$wsRA.Range("D11").Value = "Projector"
$wsRA.Range("E11").Value = 2
$wsRA.Range("F11").Value = "A"
$wsRA.Range("G11").Value = 1
$wsRA.Range("H11").Value = "EC161"
$wsRA.Range("I11").Value = "0.4"

$wsRA.Range("A12").Value = "L408"
$wsRA.Range("B12").Value = "classroom without projector"
$wsRA.Range("C12").Value = "78"
$wsRA.Range("D12").Value = "Computers"
$wsRA.Range("E12").Value = 4
$wsRA.Range("F12").Value = "A, B"
$wsRA.Range("G12").Value = 1
$wsRA.Range("H12").Value = "EC161 (Lab)"
$wsRA.Range("I12").Value = "0.8"

# ---------------------------------------------------------------------
# LTPSC_Compliance: swap the emoji glyphs for plain ASCII tags
# ---------------------------------------------------------------------
$wsLC = $wb.Worksheets.Item("LTPSC_Compliance")

$wsLC.Range("G2").Value = "[OK]"
$wsLC.Range("H2").Value = "[OK]"
$wsLC.Range("I2").Value = "[OK]"
$wsLC.Range("J2").Value = "[OK] FULLY COMPLIANT"

$wsLC.Range("G3").Value = "[FAIL]"
$wsLC.Range("H3").Value = "[OK]"
$wsLC.Range("I3").Value = "[OK]"
$wsLC.Range("J3").Value = "[WARN] PARTIAL"

$wsLC.Range("G4").Value = "[FAIL]"
$wsLC.Range("H4").Value = "[OK]"
$wsLC.Range("I4").Value = "[OK]"
$wsLC.Range("J4").Value = "[WARN] PARTIAL"

$wsLC.Range("G5").Value = "[FAIL]"
$wsLC.Range("H5").Value = "[OK]"
$wsLC.Range("I5").Value = "[OK]"
$wsLC.Range("J5").Value = "[WARN] PARTIAL"

$wsLC.Range("G6").Value = "[FAIL]"
$wsLC.Range("H6").Value = "[OK]"
$wsLC.Range("I6").Value = "[FAIL]"
$wsLC.Range("J6").Value = "[WARN] PARTIAL"

$wsLC.Range("G7").Value = "[OK]"
$wsLC.Range("H7").Value = "[OK]"
$wsLC.Range("I7").Value = "[OK]"
$wsLC.Range("J7").Value = "[OK] FULLY COMPLIANT"

# ---------------------------------------------------------------------
# Executive_Summary: refresh generation timestamp + room utilisation
# ---------------------------------------------------------------------
$wsES = $wb.Worksheets.Item("Executive_Summary")

$wsES.Range("C3").Value = "2025-12-12 16:58"
$wsES.Range("C7").Value = "11/35"
$wsES.Range("D7").Value = "Utilization: 31.4%"
$wsES.Range("C9").Value = "[WARN] NEEDS REVIEW"
